$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141, shifting the existing rows 141-174
# down to 142-175 (grows the data range from A1:R174 to A1:R175).
$ws.Rows("141:141").Insert()

# Populate the newly inserted row 141 with the new weekly record.
$ws.Range("A141").Value = 5
$ws.Range("B141").Value = "Macroferia Regional de Talca"
$ws.Range("C141").Value = "Maule"
$ws.Range("D141").Value = 44951
$ws.Range("E141").Value = 7
$ws.Range("F141").Value = 100112030
$ws.Range("G141").Value = "Poroto granado"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 100
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("M141").Value = 40000
$ws.Range("N141").Value = "`$/saco 25 kilos"
$ws.Range("O141").Value = "Región del Maule"
$ws.Range("P141").Value = 1600
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"
